# Scheduled runner update: refresh market-board derived columns (H-N)
# in each class's Leve Profits table. Values below mirror the latest
# Universalis price pull; set() for changed/added cells, clear for removed ones.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I6").Value = 700
$ws.Range("K6").Value = 2100
$ws.Range("M6").Value = -1988
$ws.Range("H62").Value = 22733256
$ws.Range("I62").Value = 1460
$ws.Range("J62").Value = 62513900
$ws.Range("K62").Value = 1460
$ws.Range("L62").Value = 62513900
$ws.Range("M62").Value = -836
$ws.Range("N62").Value = -62515148
$ws.Range("H65").Value = 22733256
$ws.Range("I65").Value = 1460
$ws.Range("J65").Value = 62513900
$ws.Range("K65").Value = 7300
$ws.Range("L65").Value = 312569500
$ws.Range("M65").Value = -4180
$ws.Range("N65").Value = -312575740
$ws.Range("H93").Value = 42500
$ws.Range("J93").Value = 42500
$ws.Range("L93").Value = 42500
$ws.Range("N93").Value = -47492
$ws.Range("H116").Value = 5172.706
$ws.Range("I116").Value = 3120
$ws.Range("K116").Value = 3120
$ws.Range("M116").Value = 322
$ws.Range("H137").Value = 387175.34
$ws.Range("I137").Value = 734105.3
$ws.Range("J137").Value = 3726.4736
$ws.Range("K137").Value = 2202315.9
$ws.Range("L137").Value = 11179.4208
$ws.Range("M137").Value = -2199765.9
$ws.Range("N137").Value = -16279.4208
$ws.Range("H141").Value = 2257.04
$ws.Range("I141").Value = 1841.4615
$ws.Range("J141").Value = 2477.551
$ws.Range("K141").Value = 5524.3845
$ws.Range("L141").Value = 7432.653
$ws.Range("M141").Value = -344.3845000000001
$ws.Range("N141").Value = -17792.653

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 22400
$ws.Range("J16").Value = 29533.334
$ws.Range("L16").Value = 29533.334
$ws.Range("N16").Value = -30107.334
$ws.Range("H32").Value = 20700.19
$ws.Range("I32").Value = 5816.0977
$ws.Range("J32").Value = 114584.46
$ws.Range("K32").Value = 5816.0977
$ws.Range("L32").Value = 114584.46
$ws.Range("M32").Value = -5529.0977
$ws.Range("N32").Value = -115158.46
$ws.Range("H44").Value = 30866
$ws.Range("J44").Value = 30866
$ws.Range("L44").Value = 30866
$ws.Range("N44").Value = -31842
$ws.Range("H63").Value = 3984.6155
$ws.Range("I63").Value = 2760
$ws.Range("J63").Value = 4750
$ws.Range("K63").Value = 2760
$ws.Range("L63").Value = 4750
$ws.Range("M63").Value = -2074
$ws.Range("N63").Value = -6122
$ws.Range("H66").Value = 3984.6155
$ws.Range("I66").Value = 2760
$ws.Range("J66").Value = 4750
$ws.Range("K66").Value = 13800
$ws.Range("L66").Value = 23750
$ws.Range("M66").Value = -10368
$ws.Range("N66").Value = -30614
$ws.Range("H74").Value = 2474.3845
$ws.Range("I74").Value = 1388.8695
$ws.Range("J74").Value = 10796.667
$ws.Range("K74").Value = 1388.8695
$ws.Range("L74").Value = 10796.667
$ws.Range("M74").Value = -514.8695
$ws.Range("N74").Value = -12544.667
$ws.Range("H77").Value = 2474.3845
$ws.Range("I77").Value = 1388.8695
$ws.Range("J77").Value = 10796.667
$ws.Range("K77").Value = 6944.3475
$ws.Range("L77").Value = 53983.335
$ws.Range("M77").Value = -2576.3475
$ws.Range("N77").Value = -62719.335
$ws.Range("H80").Value = 32304
$ws.Range("J80").Value = 38105
$ws.Range("L80").Value = 38105
$ws.Range("N80").Value = -40101
$ws.Range("H83").Value = 32304
$ws.Range("J83").Value = 38105
$ws.Range("L83").Value = 114315
$ws.Range("N83").Value = -124299
$ws.Range("H106").Value = 39888
$ws.Range("J106").Value = 39888
$ws.Range("L106").Value = 39888
$ws.Range("N106").Value = -42412

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 34271.332
$ws.Range("J82").Value = 35327.6
$ws.Range("L82").Value = 35327.6
$ws.Range("N82").Value = -36093.6
$ws.Range("H85").Value = 34271.332
$ws.Range("J85").Value = 35327.6
$ws.Range("L85").Value = 35327.6
$ws.Range("N85").Value = -37979.6
$ws.Range("H88").Value = 25343
$ws.Range("J88").Value = 25343
$ws.Range("L88").Value = 25343
$ws.Range("N88").Value = -26155
$ws.Range("H91").Value = 25343
$ws.Range("J91").Value = 25343
$ws.Range("L91").Value = 25343
$ws.Range("N91").Value = -28151
$ws.Range("H135").Value = 61638.535
$ws.Range("J135").Value = 61638.535
$ws.Range("L135").Value = 61638.535
$ws.Range("N135").Value = -71778.535

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 24900
$ws.Range("J29").Value = 24900
$ws.Range("L29").Value = 24900
$ws.Range("N29").Value = -25486
$ws.Range("H41").Value = 23960
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 23960
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 23960
$ws.Range("N41").Value = -24816
$ws.Range("M41").ClearContents()
$ws.Range("H59").Value = 16295.75
$ws.Range("J59").Value = 16295.75
$ws.Range("L59").Value = 16295.75
$ws.Range("N59").Value = -18585.75
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H109").Value = 13874.444
$ws.Range("J109").Value = 13874.444
$ws.Range("L109").Value = 13874.444
$ws.Range("N109").Value = -15954.444
$ws.Range("H134").Value = 3226.2122
$ws.Range("I134").Value = 3260.4614
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 9781.3842
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -7246.3842
$ws.Range("N134").Value = -8070

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 2576.3
$ws.Range("I119").Value = 1980.4286
$ws.Range("J119").Value = 3966.6667
$ws.Range("K119").Value = 5941.2858
$ws.Range("L119").Value = 11900.0001
$ws.Range("M119").Value = -1103.2858
$ws.Range("N119").Value = -21576.0001
$ws.Range("H129").Value = 1470.1154
$ws.Range("I129").Value = 593.0769
$ws.Range("J129").Value = 2347.1538
$ws.Range("K129").Value = 1779.2307
$ws.Range("L129").Value = 7041.4614
$ws.Range("M129").Value = 3220.7693
$ws.Range("N129").Value = -17041.4614
$ws.Range("H131").Value = 877.72
$ws.Range("I131").Value = 488.83334
$ws.Range("J131").Value = 930.75
$ws.Range("K131").Value = 1466.50002
$ws.Range("L131").Value = 2792.25
$ws.Range("M131").Value = 3573.49998
$ws.Range("N131").Value = -12872.25
$ws.Range("H140").Value = 28317.615
$ws.Range("I140").Value = 92118.09
$ws.Range("J140").Value = 3253.1428
$ws.Range("K140").Value = 276354.27
$ws.Range("L140").Value = 9759.428400000001
$ws.Range("M140").Value = -271174.27
$ws.Range("N140").Value = -20119.4284

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H123").Value = 32973.25
$ws.Range("J123").Value = 32973.25
$ws.Range("L123").Value = 32973.25
$ws.Range("N123").Value = -37873.25
$ws.Range("H140").Value = 39870.3
$ws.Range("J140").Value = 39870.3
$ws.Range("L140").Value = 39870.3
$ws.Range("N140").Value = -50230.3

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1000000000
$ws.Range("J43").Value = 1000000000
$ws.Range("L43").Value = 1000000000
$ws.Range("N43").Value = -1000000386

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1641.4736
$ws.Range("I132").Value = 1734.9025
$ws.Range("J132").Value = 1402.0625
$ws.Range("K132").Value = 5204.7075
$ws.Range("L132").Value = 4206.1875
$ws.Range("M132").Value = -2674.7075
$ws.Range("N132").Value = -9266.1875
